# "last uapdate at all"
# - Feuil1!C9 was corrected from 12 to 24.
# - A new reference row was appended at the bottom of the table (row 34):
#     A34 = "3M BPEO T2 (BDP)", B34 = 14, C34 = 24
# - The view was scrolled/selected so the next empty row (A35) is active.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing value -------------------------------------------------
$ws.Range("C9").Value = 24

# --- Append the new row --------------------------------------------------
# Clone formatting from existing rows so the new row matches the look of
# the rest of the table (thin border all around, like every data row).
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A34").PasteSpecial(-4122) | Out-Null          # xlPasteFormats

$ws.Range("B33:C33").Copy() | Out-Null
$ws.Range("B34:C34").PasteSpecial(-4122) | Out-Null      # xlPasteFormats

$excel.CutCopyMode = 0

$ws.Range("A34").Value = "3M BPEO T2 (BDP)"
$ws.Range("B34").Value = 14
$ws.Range("C34").Value = 24

# --- Update view/selection to the row right below the new data ----------
$ws.Range("A35").Select() | Out-Null
